# Updating slides on Docker and adding to index slide
#
# docker-dev-and-prod.pptx: bump the PHP version shown on every
# "Apache webserver + PHP 8.x" label (dev boxes + staging/production box)
# from 8.1 to 8.2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Apache webserver + PHP 8.1") {
            $tr.Text = "Apache webserver + PHP 8.2"
        }
    }
}
